$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting from the row above so the new row inherits the same
# date / boolean cell styles (keeps existing style indices, avoids
# creating new numFmt/style entries).
$ws.Range("A6:I6").Copy($ws.Range("A7:I7"))

$ws.Range("A7").Value = 42649.656180555554
$ws.Range("B7").Value = $false
$ws.Range("C7").Value = 9965.3799999999992
$ws.Range("D7").Value = 9994.36
$ws.Range("E7").Value = 104.82
$ws.Range("F7").Value = 105.43
$ws.Range("G7").Value = $true
$ws.Range("H7").Value = 0.57999999999999996
$ws.Range("I7").Value = $false
